$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirements")

# Copy the formatting of the last existing data row (row 43) down to the two
# new rows, then fill in their values (TC_41 / TC_42) in column C only,
# matching the pattern already used by rows 42-43 (A/C populated, B empty).
$ws.Range("C43").Copy()
$ws.Range("C44:C45").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C44").Value = "TC_41"
$ws.Range("C45").Value = "TC_42"

$ws.Range("C43:C45").Select()
